$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.735.26'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '3.464.01'
$ws.Range('E3').Value = '  -0.70%  '
$ws.Range('E4').Value = '  +0.08%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '586.70'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.44%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '178.19'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +1.96%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.628'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +5.72%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').Value = '3.463.31'
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('E10').Value = '  +1.19%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '6.98'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +1.27%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.419'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('D13').Value = '4.063.77'
$ws.Range('E13').Value = '  -0.75%  '
$ws.Range('E14').Value = '  +1.53%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '30.12'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('D16').Value = '66.539.30'
$ws.Range('E16').Value = '  +0.27%  '
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('D18').Value = '3.501.57'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('E19').Value = '  +0.02%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '13.90'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +0.47%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '371.85'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -2.07%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '7.67'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -2.00%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '73.48'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +1.98%  '
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.0000127'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +5.60%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.536'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -2.05%  '
$ws.Range('E27').Value = '  +1.54%  '
$ws.Range('E28').Value = '  +2.90%  '
$ws.Range('E29').Value = '  +0.00%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '5.93'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('E31').Value = '  +0.37%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '23.67'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -2.65%  '
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range('E35').Value = '  -3.46%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.57'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -0.28%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '162.22'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +1.69%  '
$ws.Range('E38').Value = '  -0.82%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '27.95'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -5.21%  '
$ws.Range('E40').Value = '  +2.21%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '4.51'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('D42').Value = '2.772.43'
$ws.Range('E42').Value = '  +3.56%  '
$ws.Range('E43').Value = '  +1.27%  '
$ws.Range('E44').Value = '  +0.96%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0696'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.05%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '25.52'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +4.21%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '341.39'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +8.35%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '40.05'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -1.36%  '
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('E50').Value = '  +2.91%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '31.77'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +3.17%  '
